# "32 33 leetcode done"
# Mark LeetCode problems #31, #32, #33 as solved & checked, and update the
# self-assessment note on problem #16 from "50/50, можно лучше" to
# "можно лучше".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Problems 31, 32, 33 -> mark as Решено (Solved) / Проверено (Checked)
$ws.Range("G31").Value = "Решено"
$ws.Range("H31").Value = "Проверено"

$ws.Range("G32").Value = "Решено"
$ws.Range("H32").Value = "Проверено"

$ws.Range("G33").Value = "Решено"
$ws.Range("H33").Value = "Проверено"

# Update self-assessment note on row 16 (problem #15)
$ws.Range("I16").Value = "можно лучше"

# View state left over from the editing session (zoom + active cell)
$excel.ActiveWindow.Zoom = 150
$ws.Range("H33").Select()
